$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G10").Value = "Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg"
$ws.Range("G15").Value = "nancy.abdelshafy@med.asu.edu.eg, hana.amr@med.asu.edu.eg"
$ws.Range("G18").Value = "aya.hanafy@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G19").Value = "wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Range("G20").Value = "Wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G24").Value = "servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G25").Value = "alshimaa.atef@med.asu.edu.egm, manar.montaser@med.asu.edu.eg, backup@backdoor.com"
$ws.Range("G32").Value = "Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg"
$ws.Range("G37").Value = "nancy.abdelshafy@med.asu.edu.eg, hana.amr@med.asu.edu.eg"
$ws.Range("G40").Value = "aya.hanafy@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G41").Value = "wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Range("G42").Value = "Wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G46").Value = "nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G54").Value = "Eman.m.abosakaya@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, maimustafa@med.asu.edu.eg"
$ws.Range("G62").Value = "aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Range("G63").Value = "Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Range("G64").Value = "wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"
$ws.Range("G68").Value = "nourhan.mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G76").Value = "Eman.m.abosakaya@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, maimustafa@med.asu.edu.eg"
$ws.Range("G81").Value = "Walaa.h.ghanima@med.asu.edu.eg, user@user.com, enas.omran@med.asu.edu.eg"
$ws.Range("G84").Value = "aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Range("G85").Value = "Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Range("G86").Value = "wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"
$ws.Range("G90").Value = "shaimaa.ahmed@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G96").Value = "norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G98").Value = "Eman.m.abosakaya@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, maimustafa@med.asu.edu.eg"
$ws.Range("G106").Value = "wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G107").Value = "wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G108").Value = "Wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G112").Value = "shaimaa.ahmed@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G118").Value = "norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg"
$ws.Range("G120").Value = "Eman.m.abosakaya@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, maimustafa@med.asu.edu.eg"
$ws.Range("G128").Value = "wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G129").Value = "wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
$ws.Range("G130").Value = "Wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G134").Value = "asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G142").Value = "yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg"
$ws.Range("G150").Value = "wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Range("G151").Value = "wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
$ws.Range("G156").Value = "majorelle.magdy@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, manar.montaser@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G164").Value = "yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, merna.said@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg"
$ws.Range("G172").Value = "wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
$ws.Range("G173").Value = "wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
